$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift quarter header labels left by one (drop oldest quarter, append newest)
$ws.Range("E8").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F8").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G8").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H8").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I8").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J8").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K8").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L8").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M8").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N8").Value = "فصل اول منتهی به 1401/12"
$ws.Range("E19").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F19").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G19").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H19").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I19").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J19").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K19").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L19").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M19").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N19").Value = "فصل اول منتهی به 1401/12"
$ws.Range("E31").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F31").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G31").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H31").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I31").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J31").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K31").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L31").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M31").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N31").Value = "فصل اول منتهی به 1401/12"
$ws.Range("E43").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F43").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G43").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H43").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I43").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J43").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K43").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L43").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M43").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N43").Value = "فصل اول منتهی به 1401/12"
$ws.Range("E54").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F54").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G54").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H54").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I54").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J54").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K54").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L54").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M54").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N54").Value = "فصل اول منتهی به 1401/12"
$ws.Range("E66").Value = "فصل چهارم منتهی به 1399/09"
$ws.Range("F66").Value = "فصل اول منتهی به 1399/12"
$ws.Range("G66").Value = "فصل دوم منتهی به 1400/03"
$ws.Range("H66").Value = "فصل سوم منتهی به 1400/06"
$ws.Range("I66").Value = "فصل چهارم منتهی به 1400/09"
$ws.Range("J66").Value = "فصل اول منتهی به 1400/12"
$ws.Range("K66").Value = "فصل دوم منتهی به 1401/03"
$ws.Range("L66").Value = "فصل سوم منتهی به 1401/06"
$ws.Range("M66").Value = "فصل چهارم منتهی به 1401/09"
$ws.Range("N66").Value = "فصل اول منتهی به 1401/12"
# Shift data columns E:N left by one and populate newest quarter (N) with new figures
$ws.Range("E10").Value = 13070
$ws.Range("F10").Value = 8901
$ws.Range("G10").Value = 9096
$ws.Range("H10").Value = 9815
$ws.Range("I10").Value = 14045
$ws.Range("J10").Value = 5685
$ws.Range("K10").Value = "-"
$ws.Range("L10").Value = 5876
$ws.Range("M10").Value = 2668
$ws.Range("N10").Value = 3689
$ws.Range("E11").Value = 18883
$ws.Range("F11").Value = 15311
$ws.Range("G11").Value = 13062
$ws.Range("H11").Value = 18799
$ws.Range("I11").Value = 21236
$ws.Range("J11").Value = 10578
$ws.Range("K11").Value = 13035
$ws.Range("L11").Value = 20131
$ws.Range("M11").Value = 16714
$ws.Range("N11").Value = 13277
$ws.Range("E12").Value = 6893
$ws.Range("F12").Value = 6466
$ws.Range("G12").Value = 4475
$ws.Range("H12").Value = 6700
$ws.Range("I12").Value = 8480
$ws.Range("J12").Value = 3637
$ws.Range("K12").Value = 4373
$ws.Range("L12").Value = 4727
$ws.Range("M12").Value = 5497
$ws.Range("N12").Value = 2698
$ws.Range("E13").Value = "-"
$ws.Range("F13").Value = "-"
$ws.Range("G13").Value = "-"
$ws.Range("H13").Value = 0
$ws.Range("I13").Value = "-"
$ws.Range("J13").Value = "-"
$ws.Range("K13").Value = "-"
$ws.Range("L13").Value = 0
$ws.Range("M13").Value = 0
$ws.Range("N13").Value = 0
$ws.Range("E14").Value = "-"
$ws.Range("F14").Value = "-"
$ws.Range("G14").Value = "-"
$ws.Range("H14").Value = "-"
$ws.Range("I14").Value = "-"
$ws.Range("J14").Value = "-"
$ws.Range("K14").Value = "-"
$ws.Range("L14").Value = "-"
$ws.Range("M14").Value = 2759
$ws.Range("N14").Value = 2301
$ws.Range("E15").Value = 38846
$ws.Range("F15").Value = 30678
$ws.Range("G15").Value = 26633
$ws.Range("H15").Value = 35314
$ws.Range("I15").Value = 43761
$ws.Range("J15").Value = 19900
$ws.Range("K15").Value = 17408
$ws.Range("L15").Value = 30734
$ws.Range("M15").Value = 27638
$ws.Range("N15").Value = 21965
$ws.Range("E21").Value = 6554
$ws.Range("F21").Value = 7768
$ws.Range("G21").Value = 7649
$ws.Range("H21").Value = 13342
$ws.Range("I21").Value = 12728
$ws.Range("J21").Value = 6973
$ws.Range("K21").Value = "-"
$ws.Range("L21").Value = 6047
$ws.Range("M21").Value = 2911
$ws.Range("N21").Value = 3716
$ws.Range("E22").Value = 19914
$ws.Range("F22").Value = 15003
$ws.Range("G22").Value = 11931
$ws.Range("H22").Value = 18373
$ws.Range("I22").Value = 21873
$ws.Range("J22").Value = 9396
$ws.Range("K22").Value = 14658
$ws.Range("L22").Value = 20273
$ws.Range("M22").Value = 17120
$ws.Range("N22").Value = 13898
$ws.Range("E23").Value = 7094
$ws.Range("F23").Value = 6359
$ws.Range("G23").Value = 4242
$ws.Range("H23").Value = 6604
$ws.Range("I23").Value = 7672
$ws.Range("J23").Value = 3693
$ws.Range("K23").Value = 4093
$ws.Range("L23").Value = 4984
$ws.Range("M23").Value = 4564
$ws.Range("N23").Value = 4661
$ws.Range("E24").Value = 0
$ws.Range("F24").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("H24").Value = "-"
$ws.Range("I24").Value = "-"
$ws.Range("J24").Value = "-"
$ws.Range("K24").Value = "-"
$ws.Range("L24").Value = "-"
$ws.Range("M24").Value = "-"
$ws.Range("N24").Value = "-"
$ws.Range("E25").Value = 0
$ws.Range("F25").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("H25").Value = 0
$ws.Range("I25").Value = "-"
$ws.Range("J25").Value = "-"
$ws.Range("K25").Value = "-"
$ws.Range("L25").Value = 0
$ws.Range("M25").Value = 0
$ws.Range("N25").Value = 0
$ws.Range("E26").Value = "-"
$ws.Range("F26").Value = "-"
$ws.Range("G26").Value = "-"
$ws.Range("H26").Value = "-"
$ws.Range("I26").Value = "-"
$ws.Range("J26").Value = "-"
$ws.Range("K26").Value = "-"
$ws.Range("L26").Value = "-"
$ws.Range("M26").Value = 2759
$ws.Range("N26").Value = 2301
$ws.Range("E27").Value = 33562
$ws.Range("F27").Value = 29130
$ws.Range("G27").Value = 23822
$ws.Range("H27").Value = 38319
$ws.Range("I27").Value = 42273
$ws.Range("J27").Value = 20062
$ws.Range("K27").Value = 18751
$ws.Range("L27").Value = 31304
$ws.Range("M27").Value = 27354
$ws.Range("N27").Value = 24576
$ws.Range("E33").Value = 1879234
$ws.Range("F33").Value = 1829607
$ws.Range("G33").Value = 1847175
$ws.Range("H33").Value = 2283694
$ws.Range("I33").Value = 3173404
$ws.Range("J33").Value = 1822183
$ws.Range("K33").Value = "-"
$ws.Range("L33").Value = 2170893
$ws.Range("M33").Value = 1394719
$ws.Range("N33").Value = 1276383
$ws.Range("E34").Value = 2096076
$ws.Range("F34").Value = 1603362
$ws.Range("G34").Value = 1124022
$ws.Range("H34").Value = 2435605
$ws.Range("I34").Value = 2832947
$ws.Range("J34").Value = 1199009
$ws.Range("K34").Value = 1653799
$ws.Range("L34").Value = 2567211
$ws.Range("M34").Value = 2656497
$ws.Range("N34").Value = 2485794
$ws.Range("E35").Value = 1814663
$ws.Range("F35").Value = 1822146
$ws.Range("G35").Value = 1225489
$ws.Range("H35").Value = 2232282
$ws.Range("I35").Value = 2646316
$ws.Range("J35").Value = 1335934
$ws.Range("K35").Value = "-"
$ws.Range("L35").Value = 2201313
$ws.Range("M35").Value = 2171816
$ws.Range("N35").Value = 2466128
$ws.Range("E36").Value = 0
$ws.Range("F36").Value = 0
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = "-"
$ws.Range("I36").Value = "-"
$ws.Range("J36").Value = "-"
$ws.Range("K36").Value = "-"
$ws.Range("L36").Value = "-"
$ws.Range("M36").Value = "-"
$ws.Range("N36").Value = "-"
$ws.Range("E37").Value = 0
$ws.Range("F37").Value = 0
$ws.Range("G37").Value = 0
$ws.Range("H37").Value = 0
$ws.Range("I37").Value = "-"
$ws.Range("J37").Value = "-"
$ws.Range("K37").Value = "-"
$ws.Range("L37").Value = 0
$ws.Range("M37").Value = 0
$ws.Range("N37").Value = 0
$ws.Range("E38").Value = "-"
$ws.Range("F38").Value = "-"
$ws.Range("G38").Value = "-"
$ws.Range("H38").Value = "-"
$ws.Range("I38").Value = "-"
$ws.Range("J38").Value = "-"
$ws.Range("K38").Value = "-"
$ws.Range("L38").Value = "-"
$ws.Range("M38").Value = 130621
$ws.Range("N38").Value = 129789
$ws.Range("E39").Value = 5789973
$ws.Range("F39").Value = 5255115
$ws.Range("G39").Value = 4196686
$ws.Range("H39").Value = 6951581
$ws.Range("I39").Value = 8652667
$ws.Range("J39").Value = 4357126
$ws.Range("K39").Value = 1653799
$ws.Range("L39").Value = 6939417
$ws.Range("M39").Value = 6353653
$ws.Range("N39").Value = 6358094
$ws.Range("E45").Value = 156290253
$ws.Range("F45").Value = 235531282
$ws.Range("G45").Value = 241492352
$ws.Range("H45").Value = 171165792
$ws.Range("I45").Value = 249324639
$ws.Range("J45").Value = 261319805
$ws.Range("K45").Value = "-"
$ws.Range("L45").Value = 359003307
$ws.Range("M45").Value = 479120234
$ws.Range("N45").Value = 343483046
$ws.Range("E46").Value = 105256403
$ws.Range("F46").Value = 106869426
$ws.Range("G46").Value = 94210209
$ws.Range("H46").Value = 132564361
$ws.Range("I46").Value = 129517990
$ws.Range("J46").Value = 127608450
$ws.Range("K46").Value = 112825692
$ws.Range("L46").Value = 126632023
$ws.Range("M46").Value = 155169217
$ws.Range("N46").Value = 178859836
$ws.Range("E47").Value = 255802509
$ws.Range("F47").Value = 286545998
$ws.Range("G47").Value = 288894154
$ws.Range("H47").Value = 338019685
$ws.Range("I47").Value = 344931700
$ws.Range("J47").Value = 361747631
$ws.Range("K47").Value = 383192768
$ws.Range("L47").Value = 441675963
$ws.Range("M47").Value = 475858019
$ws.Range("N47").Value = 529098477
$ws.Range("E48").Value = "-"
$ws.Range("F48").Value = "-"
$ws.Range("G48").Value = "-"
$ws.Range("H48").Value = "-"
$ws.Range("I48").Value = "-"
$ws.Range("J48").Value = "-"
$ws.Range("K48").Value = "-"
$ws.Range("L48").Value = "-"
$ws.Range("M48").Value = "-"
$ws.Range("N48").Value = "-"
$ws.Range("E49").Value = "-"
$ws.Range("F49").Value = "-"
$ws.Range("G49").Value = "-"
$ws.Range("H49").Value = "-"
$ws.Range("I49").Value = "-"
$ws.Range("J49").Value = "-"
$ws.Range("K49").Value = "-"
$ws.Range("L49").Value = "-"
$ws.Range("M49").Value = "-"
$ws.Range("N49").Value = "-"
$ws.Range("E50").Value = "-"
$ws.Range("F50").Value = "-"
$ws.Range("G50").Value = "-"
$ws.Range("H50").Value = "-"
$ws.Range("I50").Value = "-"
$ws.Range("J50").Value = "-"
$ws.Range("K50").Value = "-"
$ws.Range("L50").Value = "-"
$ws.Range("M50").Value = 47343603
$ws.Range("N50").Value = 56405476
$ws.Range("E56").Value = -1354680
$ws.Range("F56").Value = -1528202
$ws.Range("G56").Value = -1593862
$ws.Range("H56").Value = -1882002
$ws.Range("I56").Value = -2561059
$ws.Range("J56").Value = -1612458
$ws.Range("K56").Value = "-"
$ws.Range("L56").Value = -1725169
$ws.Range("M56").Value = -1248730
$ws.Range("N56").Value = -1093779
$ws.Range("E57").Value = -1571799
$ws.Range("F57").Value = -1378362
$ws.Range("G57").Value = -1062754
$ws.Range("H57").Value = -2133453
$ws.Range("I57").Value = -2768399
$ws.Range("J57").Value = -1233163
$ws.Range("K57").Value = -1910780
$ws.Range("L57").Value = -2614953
$ws.Range("M57").Value = -2538572
$ws.Range("N57").Value = -2253009
$ws.Range("E58").Value = -1452652
$ws.Range("F58").Value = -1674727
$ws.Range("G58").Value = -1137089
$ws.Range("H58").Value = -1793681
$ws.Range("I58").Value = -2383859
$ws.Range("J58").Value = -1329646
$ws.Range("K58").Value = "-"
$ws.Range("L58").Value = -2192149
$ws.Range("M58").Value = -1994107
$ws.Range("N58").Value = -2434485
$ws.Range("E59").Value = 0
$ws.Range("F59").Value = 0
$ws.Range("G59").Value = 0
$ws.Range("H59").Value = "-"
$ws.Range("I59").Value = "-"
$ws.Range("J59").Value = "-"
$ws.Range("K59").Value = "-"
$ws.Range("L59").Value = "-"
$ws.Range("M59").Value = "-"
$ws.Range("N59").Value = "-"
$ws.Range("E60").Value = 0
$ws.Range("F60").Value = 0
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 0
$ws.Range("I60").Value = "-"
$ws.Range("J60").Value = "-"
$ws.Range("K60").Value = "-"
$ws.Range("L60").Value = 0
$ws.Range("M60").Value = 0
$ws.Range("N60").Value = 0
$ws.Range("E61").Value = "-"
$ws.Range("F61").Value = "-"
$ws.Range("G61").Value = "-"
$ws.Range("H61").Value = "-"
$ws.Range("I61").Value = "-"
$ws.Range("J61").Value = "-"
$ws.Range("K61").Value = "-"
$ws.Range("L61").Value = "-"
$ws.Range("M61").Value = -67424
$ws.Range("N61").Value = -84044
$ws.Range("E62").Value = -4379131
$ws.Range("F62").Value = -4581291
$ws.Range("G62").Value = -3793705
$ws.Range("H62").Value = -5809136
$ws.Range("I62").Value = -7713317
$ws.Range("J62").Value = -4175267
$ws.Range("K62").Value = -1910780
$ws.Range("L62").Value = -6532271
$ws.Range("M62").Value = -5848833
$ws.Range("N62").Value = -5865317
$ws.Range("E68").Value = 524554
$ws.Range("F68").Value = 301405
$ws.Range("G68").Value = 253313
$ws.Range("H68").Value = 401692
$ws.Range("I68").Value = 612345
$ws.Range("J68").Value = 209725
$ws.Range("K68").Value = "-"
$ws.Range("L68").Value = 445724
$ws.Range("M68").Value = 145989
$ws.Range("N68").Value = 182604
$ws.Range("E69").Value = 524277
$ws.Range("F69").Value = 225000
$ws.Range("G69").Value = 61268
$ws.Range("H69").Value = 302152
$ws.Range("I69").Value = 64548
$ws.Range("J69").Value = -34154
$ws.Range("K69").Value = -256981
$ws.Range("L69").Value = -47742
$ws.Range("M69").Value = 117925
$ws.Range("N69").Value = 232785
$ws.Range("E70").Value = 362011
$ws.Range("F70").Value = 147419
$ws.Range("G70").Value = 88400
$ws.Range("H70").Value = 438601
$ws.Range("I70").Value = 262457
$ws.Range("J70").Value = 6288
$ws.Range("K70").Value = "-"
$ws.Range("L70").Value = 9164
$ws.Range("M70").Value = 177709
$ws.Range("N70").Value = 31643
$ws.Range("E71").Value = 0
$ws.Range("F71").Value = 0
$ws.Range("G71").Value = 0
$ws.Range("H71").Value = "-"
$ws.Range("I71").Value = "-"
$ws.Range("J71").Value = "-"
$ws.Range("K71").Value = "-"
$ws.Range("L71").Value = "-"
$ws.Range("M71").Value = "-"
$ws.Range("N71").Value = "-"
$ws.Range("E72").Value = 0
$ws.Range("F72").Value = 0
$ws.Range("G72").Value = 0
$ws.Range("H72").Value = 0
$ws.Range("I72").Value = "-"
$ws.Range("J72").Value = "-"
$ws.Range("K72").Value = "-"
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = 0
$ws.Range("N72").Value = 0
$ws.Range("E73").Value = "-"
$ws.Range("F73").Value = "-"
$ws.Range("G73").Value = "-"
$ws.Range("H73").Value = "-"
$ws.Range("I73").Value = "-"
$ws.Range("J73").Value = "-"
$ws.Range("K73").Value = "-"
$ws.Range("L73").Value = "-"
$ws.Range("M73").Value = 63197
$ws.Range("N73").Value = 45745
$ws.Range("E74").Value = 1410842
$ws.Range("F74").Value = 673824
$ws.Range("G74").Value = 402981
$ws.Range("H74").Value = 1142445
$ws.Range("I74").Value = 939350
$ws.Range("J74").Value = 181859
$ws.Range("K74").Value = -256981
$ws.Range("L74").Value = 407146
$ws.Range("M74").Value = 504820
$ws.Range("N74").Value = 492777